$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 279.8387
$ws.Range("I28").Value = 264.92
$ws.Range("J28").Value = 342
$ws.Range("K28").Value = 264.92
$ws.Range("L28").Value = 342
$ws.Range("M28").Value = 220.08
$ws.Range("N28").Value = -1312

$ws.Range("H61").Value = 936
$ws.Range("I61").Value = 420
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 1260
$ws.Range("L61").Value = 9000
$ws.Range("M61").Value = -1088
$ws.Range("N61").Value = -9344

$ws.Range("H121").Value = 2550
$ws.Range("J121").Value = 3066.6667
$ws.Range("L121").Value = 9200.000100000001
$ws.Range("N121").Value = -12694.0001

$ws.Range("H135").Value = 1203.3478
$ws.Range("I135").Value = 1026.5555
$ws.Range("J135").Value = 1839.8
$ws.Range("K135").Value = 9238.9995
$ws.Range("L135").Value = 16558.2
$ws.Range("M135").Value = -6703.9995
$ws.Range("N135").Value = -21628.2

$ws.Range("H138").Value = 2729.3223
$ws.Range("I138").Value = 1128.16
$ws.Range("J138").Value = 3345.1538
$ws.Range("K138").Value = 3384.48
$ws.Range("L138").Value = 10035.4614
$ws.Range("M138").Value = 1755.52
$ws.Range("N138").Value = -20315.4614

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 767.0952
$ws.Range("I2").Value = 576.35297
$ws.Range("K2").Value = 576.35297
$ws.Range("M2").Value = -463.35297

$ws.Range("H61").Value = 2544.818
$ws.Range("I61").Value = 2554.7778
$ws.Range("J61").Value = 2500
$ws.Range("K61").Value = 2554.7778
$ws.Range("L61").Value = 2500
$ws.Range("M61").Value = -2342.7778
$ws.Range("N61").Value = -2924

$ws.Range("H116").Value = 767.0952
$ws.Range("I116").Value = 576.35297
$ws.Range("K116").Value = 576.35297
$ws.Range("M116").Value = 1717.64703

$ws.Range("H136").Value = 2544.818
$ws.Range("I136").Value = 2554.7778
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 7664.3334
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -5114.3334
$ws.Range("N136").Value = -12600

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 767.0952
$ws.Range("I3").Value = 576.35297
$ws.Range("K3").Value = 576.35297
$ws.Range("M3").Value = -462.35297

$ws.Range("H64").Value = 234
$ws.Range("J64").Value = 246.66667
$ws.Range("L64").Value = 246.66667
$ws.Range("N64").Value = -696.6666700000001

$ws.Range("H67").Value = 234
$ws.Range("J67").Value = 246.66667
$ws.Range("L67").Value = 246.66667
$ws.Range("N67").Value = -1806.66667

$ws.Range("H134").Value = 1506.6327
$ws.Range("I134").Value = 1383.561
$ws.Range("J134").Value = 2137.375
$ws.Range("K134").Value = 4150.683
$ws.Range("L134").Value = 6412.125
$ws.Range("M134").Value = -1615.683
$ws.Range("N134").Value = -11482.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2005.4222
$ws.Range("I31").Value = 1506.2778
$ws.Range("J31").Value = 4002
$ws.Range("K31").Value = 1506.2778
$ws.Range("L31").Value = 4002
$ws.Range("M31").Value = -1211.2778
$ws.Range("N31").Value = -4592

$ws.Range("H34").Value = 2005.4222
$ws.Range("I34").Value = 1506.2778
$ws.Range("J34").Value = 4002
$ws.Range("K34").Value = 1506.2778
$ws.Range("L34").Value = 4002
$ws.Range("M34").Value = -1304.2778
$ws.Range("N34").Value = -4406

$ws.Range("H134").Value = 1870.4166
$ws.Range("I134").Value = 1717.1666
$ws.Range("J134").Value = 2330.1667
$ws.Range("K134").Value = 5151.4998
$ws.Range("L134").Value = 6990.500100000001
$ws.Range("M134").Value = -2616.4998
$ws.Range("N134").Value = -12060.5001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 5633.3335
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 5633.3335
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 16900.0005
$ws.Range("M76").Value = $null
$ws.Range("N76").Value = -17666.0005

$ws.Range("H79").Value = 5633.3335
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 5633.3335
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 16900.0005
$ws.Range("M79").Value = $null
$ws.Range("N79").Value = -19552.0005

$ws.Range("H87").Value = 2000
$ws.Range("I87").Value = 2000
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 6000
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = -4752
$ws.Range("N87").Value = $null

$ws.Range("H90").Value = 2000
$ws.Range("I90").Value = 2000
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 18000
$ws.Range("L90").Value = 0
$ws.Range("M90").Value = -11760
$ws.Range("N90").Value = $null

$ws.Range("H129").Value = 1725376.6
$ws.Range("I129").Value = 798.4286
$ws.Range("J129").Value = 2274106
$ws.Range("K129").Value = 2395.2858
$ws.Range("L129").Value = 6822318
$ws.Range("M129").Value = 2604.7142
$ws.Range("N129").Value = -6832318

$ws.Range("H137").Value = 2670.2
$ws.Range("I137").Value = 1686.6666
$ws.Range("J137").Value = 3325.889
$ws.Range("K137").Value = 5059.9998
$ws.Range("L137").Value = 9977.667000000001
$ws.Range("M137").Value = 40.0002000000004
$ws.Range("N137").Value = -20177.667

$ws.Range("H139").Value = 1608.6786
$ws.Range("I139").Value = 1021.5789
$ws.Range("J139").Value = 2848.111
$ws.Range("K139").Value = 3064.7367
$ws.Range("L139").Value = 8544.332999999999
$ws.Range("M139").Value = 2075.2633
$ws.Range("N139").Value = -18824.333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 14149.5625
$ws.Range("J123").Value = 14149.5625
$ws.Range("L123").Value = 14149.5625
$ws.Range("N123").Value = -19049.5625

$ws.Range("H126").Value = 2919.8
$ws.Range("I126").Value = 2399.75
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 7199.25
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -4729.25
$ws.Range("N126").Value = -19940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3421.8572
$ws.Range("I136").Value = 2656.9333
$ws.Range("J136").Value = 5334.1665
$ws.Range("K136").Value = 7970.7999
$ws.Range("L136").Value = 16002.4995
$ws.Range("M136").Value = -5420.7999
$ws.Range("N136").Value = -21102.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 10312.053
$ws.Range("I126").Value = 11442.883
$ws.Range("J126").Value = 700
$ws.Range("K126").Value = 34328.649
$ws.Range("L126").Value = 2100
$ws.Range("M126").Value = -31858.649
$ws.Range("N126").Value = -7040

$ws.Range("H133").Value = 29800
$ws.Range("J133").Value = 29800
$ws.Range("L133").Value = 29800
$ws.Range("N133").Value = -39920
